# Update from Github Action
#
# Insert a new data row right below the header (new row 2), pushing the
# existing company rows down by one, and fill the new row with the
# newly-added company's information ("慧资环球"). The index column (A)
# is kept as a simple sequential 0-based series across all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts old rows 2..10 down
# to 3..11).
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header); clear
# it so the new data cells start out unstyled, like the other data rows.
$ws.Range("B2:P2").ClearFormats()

# Give A2 the same style as the other index-column cells (centered/bordered
# header-like style used for column A throughout the data rows).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row 2 with the new company's data.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "慧资环球"
$ws.Range("C2").Value = "白下（年中搬到河西）"
$ws.Range("D2").Value = "研发中心"
$ws.Range("E2").Value = ".NET/Python etc."
$ws.Range("F2").Value = "自己安排，满8小时工时就好"
$ws.Range("G2").Value = "自己安排"
$ws.Range("H2").Value = "不加班"
$ws.Range("I2").Value = "全额8%"
$ws.Range("J2").Value = "13薪，每年调薪一次"
$ws.Range("K2").Value = "不打折"
$ws.Range("L2").Value = "一个高配台式机或者一个高配Dell工作站笔记本，两个40寸4K显示器 Processor Intel(R) Core(TM) i9-10980XE CPU @ 3.00GHz 3.00 GHz  128GB RAM (新的台式机配置标准)"
$ws.Range("M2").Value = "10 ~ 20天"
$ws.Range("N2").Value = "完全不打卡"
$ws.Range("O2").Value = "内推VX：Just1n"

# Renumber the index column (A) for all rows so it stays a simple sequential
# 0-based series after the insertion (0..9 across rows 2..11).
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
